# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the content is
# now "Ready for handoff" (was "In Translation"), and refreshes the
# handoff timestamps. Also widens the "Status"/locale-status columns so
# the longer "Ready for handoff" text fits.

$wb = $excel.ActiveWorkbook

$newStatus   = "Ready for handoff"
$overviewDate = "2016-09-02 04:42:25"
$zhHandoffDate = "2016-09-02 04:42:21"
$deHandoffDate = "2016-09-02 04:42:25"

# Column width used by the Status-ish columns after the content grew from
# "In Translation" to the longer "Ready for handoff". 16 + 1/3 is the
# ColumnWidth value that lands on the widened grid step.
$newColumnWidth = 16 + 1/3

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $overviewDate

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $zhHandoffDate

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $deHandoffDate

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
